$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(285,44180,160.54797487546026,89.169768851181715,22.168480451430874,138.66952215232232,103.3785546569267,98.856069660235619,133.37615795523701,28.598072738576313,138.54056504003626,119.98942242006214,0,0,107.31310034253958,14.970059880239521,114.83901671151136),
    @(286,44181,167.58717782109596,91.766946584711278,19.481391911863497,133.5831883282024,96.856248369423426,96.124295714529623,123.74061802691189,29.530618588747281,130.23115274210605,118.88759117653539,0,0,106.59448583131721,14.970059880239521,111.03179397067706),
    @(287,44182,162.71388347411738,91.766946584711278,24.183796856106408,128.49685450408245,84.789981737542391,102.95373057879461,124.24775170735006,34.193347839602112,126.15198670494033,117.23484431124528,0,0,101.08510791194577,0,109.65899730931855),
    @(288,44183,144.84513753519602,100.42420569647649,28.886201800349323,118.05648507562576,83.159405165666584,100.05122076148199,117.99310298194604,34.815045073049426,117.08717328901646,113.04788558584367,0,0,106.59448583131721,0,104.95488074973002),
    @(289,44184,138.34741173922461,87.438317028828678,26.199113260781942,100.38816758131442,75.332637620662666,94.41693699846337,103.11718168909324,34.193347839602112,112.55476658105454,108.6405606117367,0,0,87.43143219872087,0,95.949334651218123),
)

foreach ($rowData in $newRows) {
    $r = $rowData[0]
    for ($col = 1; $col -le 16; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col]
    }
}

[void]$ws.Range("A1:P289").Select()
